# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Sun Sep 24 22:43:58 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "26.619.70"
$ws.Cells.Item(2, 5).Value = "  -0.35%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.595.48"
$ws.Cells.Item(3, 5).Value = "  -0.26%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.09%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "210.59"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.39%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.511"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.27%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.0615"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.50%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.246"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.44%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "19.61"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.40%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +0.47%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "1.819.58"
$ws.Cells.Item(12, 5).Value = "  -0.27%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.587.23"
$ws.Cells.Item(13, 5).Value = "  -1.42%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +0.02%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  -0.25%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "64.53"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -1.12%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "26.589.67"

# Row 18
$ws.Cells.Item(18, 4).Value = "0.0₃0739"
$ws.Cells.Item(18, 5).Value = "  -2.05%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  -0.07%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "208.89"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.04%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "7.07"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -1.83%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  +0.16%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -3.54%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "8.95"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.18%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "144.96"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.89%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -0.07%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +0.17%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -0.48%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -0.26%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -2.77%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -0.33%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.25"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.14%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "2.96"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.22%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "1.282.45"
$ws.Cells.Item(34, 5).Value = "  -0.67%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +0.40%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +11.78%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -3.45%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.49"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.89%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -1.77%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -0.42%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.55%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "2.15"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -1.72%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.770"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -1.87%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "62.72"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.64%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "1.731.51"
$ws.Cells.Item(45, 5).Value = "  -0.36%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "89.30"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -2.10%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +0.21%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +2.60%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0513"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.56%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "7.47"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +1.14%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  -0.05%  "

Write-Output "Applied all cryptos updates"
